# Updates the cryptocurrency price table (columns D = Price, E = Volume(1h))
# on the active worksheet to reflect the latest scraped values, and swaps
# the RenderToken / BabyDogeCoin rows (48-49) which changed order.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking text (e.g. "1.000", "29.292.86") that must
# stay plain text -- prefixing with an apostrophe forces Excel to keep it as
# text instead of silently re-parsing it as a number and losing formatting.
function Set-TextCell($rangeAddr, $text) {
    $ws.Range($rangeAddr).Value = "'" + $text
}

Set-TextCell 'D2' '29.292.86'
$ws.Range('E2').Value = '  +0.26%  '
Set-TextCell 'D3' '1.872.82'
$ws.Range('E3').Value = '  +0.26%  '
Set-TextCell 'D4' '1.000'
$ws.Range('E4').Value = '  +0.06%  '
Set-TextCell 'D5' '0.7097'
$ws.Range('E5').Value = '  -0.03%  '
Set-TextCell 'D6' '241.70'
$ws.Range('E6').Value = '  -0.12%  '
$ws.Range('E7').Value = '  +0.03%  '
Set-TextCell 'D8' '0.07810'
$ws.Range('E8').Value = '  +1.98%  '
Set-TextCell 'D9' '0.3097'
$ws.Range('E9').Value = '  -0.54%  '
$ws.Range('E10').Value = '  +1.44%  '
$ws.Range('E11').Value = '  +0.28%  '
Set-TextCell 'D12' '1.871.85'
$ws.Range('E12').Value = '  +0.21%  '
Set-TextCell 'D13' '5.236'
$ws.Range('E13').Value = '  +0.07%  '
Set-TextCell 'D14' '0.7109'
$ws.Range('E14').Value = '  +0.00%  '
Set-TextCell 'D15' '91.10'
$ws.Range('E15').Value = '  -0.15%  '
Set-TextCell 'D16' '29.308.40'
$ws.Range('E16').Value = '  +0.30%  '
Set-TextCell 'D17' '6.076'
$ws.Range('E17').Value = '  +2.10%  '
Set-TextCell 'D18' '0.000008171'
$ws.Range('E18').Value = '  +4.29%  '
Set-TextCell 'D19' '239.82'
$ws.Range('E19').Value = '  -1.65%  '
$ws.Range('E20').Value = '  +0.95%  '
Set-TextCell 'D21' '2.123.60'
$ws.Range('E21').Value = '  +0.46%  '
Set-TextCell 'D22' '1.001'
$ws.Range('E22').Value = '  +0.11%  '
Set-TextCell 'D23' '7.747'
$ws.Range('E23').Value = '  -1.48%  '
Set-TextCell 'D24' '1.001'
$ws.Range('E24').Value = '  +0.06%  '
Set-TextCell 'D25' '0.1596'
$ws.Range('E25').Value = '  -2.16%  '
Set-TextCell 'D26' '162.94'
$ws.Range('E26').Value = '  -0.08%  '
Set-TextCell 'D27' '8.997'
$ws.Range('E27').Value = '  +0.36%  '
Set-TextCell 'D28' '18.45'
$ws.Range('E28').Value = '  -0.45%  '
Set-TextCell 'D29' '1.505'
$ws.Range('E29').Value = '  -0.27%  '
Set-TextCell 'D30' '4.390'
$ws.Range('E30').Value = '  -0.22%  '
Set-TextCell 'D31' '1.298'
$ws.Range('E31').Value = '  -1.39%  '
Set-TextCell 'D32' '4.295'
$ws.Range('E32').Value = '  +1.22%  '
Set-TextCell 'D33' '0.05388'
$ws.Range('E33').Value = '  +4.88%  '
$ws.Range('E34').Value = '  +1.54%  '
$ws.Range('E35').Value = '  +0.93%  '
Set-TextCell 'D36' '0.7485'
$ws.Range('E36').Value = '  -5.89%  '
Set-TextCell 'D37' '2.692'
$ws.Range('E37').Value = '  +0.18%  '
Set-TextCell 'D38' '0.01869'
$ws.Range('E38').Value = '  +0.67%  '
Set-TextCell 'D39' '1.237.01'
$ws.Range('E39').Value = '  +7.15%  '
Set-TextCell 'D40' '2.729'
$ws.Range('E40').Value = '  +0.79%  '
Set-TextCell 'D41' '6.519'
$ws.Range('E41').Value = '  +2.17%  '
Set-TextCell 'D42' '0.8882'
$ws.Range('E42').Value = '  -0.85%  '
Set-TextCell 'D43' '72.40'
$ws.Range('E43').Value = '  -1.16%  '
Set-TextCell 'D44' '108.30'
$ws.Range('E44').Value = '  +4.91%  '
Set-TextCell 'D45' '1.000'
$ws.Range('E45').Value = '  +0.05%  '
Set-TextCell 'D46' '2.019.32'
$ws.Range('E46').Value = '  +0.38%  '
$ws.Range('E47').Value = '  +0.60%  '
$ws.Range('B48').Value = 'RenderToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextCell 'D48' '1.790'
$ws.Range('E48').Value = '  +0.60%  '
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextCell 'D49' '0.00000000122'
$ws.Range('E49').Value = '  +2.41%  '
Set-TextCell 'D50' '9.410'
Set-TextCell 'D51' '0.4312'
$ws.Range('E51').Value = '  +0.39%  '

Write-Output "Applied 96 cell updates"
